# nov-17 - making generic methods for popup
# Adds a new "SelectDropndown" / "Products" column (M) to the TC04 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC04")

# New header + value cells in column M
$ws.Range("M1").Value = "SelectDropndown"
$ws.Range("M2").Value = "Products"

# Match the column width used by the other header columns on this sheet
$ws.Columns.Item(13).ColumnWidth = 16.14

# Move the active selection to the newly populated cell
$ws.Range("M2").Select()
